$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    # Force the cell to Text format first so a numeric-looking string
    # (e.g. "250.97") is not auto-converted to a Number by Excel, then
    # restore the "Normal" style so no stray number-format style sticks
    # around on the cell (keeps the style index unchanged, matching the
    # original formatting).
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

# Column D (Price) updates
Set-TextValue "D2"  "250.97"
Set-TextValue "D3"  "23.79"
Set-TextValue "D4"  "6.038"
Set-TextValue "D5"  "0.05982"
Set-TextValue "D6"  "3.430"
Set-TextValue "D7"  "6.571"
Set-TextValue "D8"  "1.317"
Set-TextValue "D9"  "0.8000"
Set-TextValue "D10" "0.1511"
Set-TextValue "D11" "0.07933"
Set-TextValue "D12" "0.03344"
Set-TextValue "D13" "0.03072"
Set-TextValue "D14" "0.09282"
Set-TextValue "D15" "3.575"
Set-TextValue "D16" "0.001657"
Set-TextValue "D17" "0.04762"
Set-TextValue "D18" "0.0006102"
Set-TextValue "D19" "0.006226"
Set-TextValue "D20" "0.005698"
Set-TextValue "D21" "0.001075"
Set-TextValue "D22" "0.0001505"
Set-TextValue "D23" "3.680"
Set-TextValue "D24" "2.204"
Set-TextValue "D25" "0.3319"
Set-TextValue "D26" "0.1232"
Set-TextValue "D27" "0.0006499"
Set-TextValue "D40" "0.04459"
Set-TextValue "D41" "0.007088"
Set-TextValue "D42" "0.1071"
Set-TextValue "D43" "0.003363"
Set-TextValue "D45" "0.002469"
Set-TextValue "D46" "0.00005908"
Set-TextValue "D47" "0.00000000753"
Set-TextValue "D48" "0.7027"
Set-TextValue "D49" "0.09558"
Set-TextValue "D50" "0.00002108"
Set-TextValue "D51" "0.01014"

# Column E (Volume(1h) label) updates
Set-TextValue "E18" "17OneONE"
Set-TextValue "E49" "48BOLOBOLOWorstin24h"
